$d = $word.ActiveDocument

# Locate the "TEAM" run inside the "TEAM VAJEF" paragraph and position
# a collapsed range right after it (before the existing " VAJEF").
$rng = $d.Content
$rng.Find.Execute("TEAM VAJEF")
$rng.Collapse(1)
$rng.Find.Execute("TEAM")
$rng.Collapse(0)

# Insert the new " -" text. It inherits the surrounding (Arial, 26pt)
# character formatting, so without the toggle below it would just merge
# back into a single "TEAM - VAJEF" run.
$rng.InsertAfter(" -")

# Toggling Bold on then back off forces the inserted text to stay in
# its own run instead of being re-merged with its neighbours, matching
# the three-run split ("TEAM" / " -" / " VAJEF") produced when Word
# itself records this kind of mid-word insertion.
$rng.Bold = 1
$rng.Bold = 0
